# Append a new data row (row 67) to the "Prices" sheet with the
# 2025-05-07 Argent (silver) price data, mirroring the previous day's
# row (row 66) for every column except the date.
#
# The source cells are all plain text (stored as inline/shared strings
# in the OOXML), even though several of them look like numbers or
# dates ("38", "5,386", "2025-05-07", ...). A direct
# `$ws.Range(...).Value = "..."` assignment would let Excel's normal
# type-inference kick in and silently convert these into real dates /
# numbers (with an auto-picked NumberFormat), which does not match the
# source data. To avoid that, we build each value as a text formula
# (="literal") - which always evaluates to a string - and then do a
# copy / paste-special-values over itself to flatten the formula down
# to a plain literal text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 67
$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rowValues = @(
    "2025-05-07",
    "38",
    "37.28",
    "1.03",
    "0.27",
    "0.09",
    "5,386",
    "8,064",
    "8,114",
    "7.2536"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = $cols[$i] + $newRow
    $escaped = $rowValues[$i].Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
}

$rowRange = $ws.Range("A" + $newRow + ":J" + $newRow)
$rowRange.Copy()
$rowRange.PasteSpecial(-4163)
